# edit.ps1 -- applies the ModelComponentClassDiagram "AddressBook -> GradTrak"
# rename pass (DG / diagram refresh) to the single slide carried in this
# deck, using PowerPoint COM-interop calls only.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# --- simple single-run label renames ---------------------------------
(Get-ShapeById $s 46).TextFrame.TextRange.Text  = "VersionedGradTrak"      # VersionedAddressBook
(Get-ShapeById $s 49).TextFrame.TextRange.Text  = "UniqueModuleTakenList"  # UniquePersonList
(Get-ShapeById $s 62).TextFrame.TextRange.Text  = "ModuleTaken"            # Person
(Get-ShapeById $s 76).TextFrame.TextRange.Text  = "Code"                   # Name
(Get-ShapeById $s 80).TextFrame.TextRange.Text  = "Semester"               # Phone
(Get-ShapeById $s 83).TextFrame.TextRange.Text  = "Grade Range"            # Email
(Get-ShapeById $s 85).TextFrame.TextRange.Text  = "Workload"               # Address
(Get-ShapeById $s 55).TextFrame.TextRange.Text  = "GradTrak"               # AddressBook

# --- shape 100 keeps its leading "<<interface>>" + line break, only the
#     second run (the class name) changes, so replace that sub-range only
$sh100 = Get-ShapeById $s 100
$tr100 = $sh100.TextFrame.TextRange
$full100 = $tr100.Text
$idx100 = $full100.IndexOf("ReadOnlyAddressBook")
$tr100.Characters($idx100 + 1, "ReadOnlyAddressBook".Length).Text = "ReadOnlyGradTrak"

# --- resize the "Workload" rectangle (id 85): cx 708186 -> 745804 EMU ----
$sh85 = Get-ShapeById $s 85
$sh85.Width = 58.72473

# --- bent connector feeding the "Workload" rectangle (id 86) now carries
#     an explicit adj1 = 50000 (50%) adjustment value
$sh86 = Get-ShapeById $s 86
$sh86.Adjustments.Item(1) = 0.5
